# Remove US Core QuestionnaireResponse Tag Element [FHIR-40742]
# Target shape: slide 2, shape index 3 ("Google Shape;191;p26")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(3)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# --- 1. Resize the shape (a:ext cx/cy) -----------------------------------
# Shape.Width/Height are expressed in points (1 pt = 12700 EMU). A tiny
# epsilon is added to counter float truncation in the EMU<->pt round trip
# so the saved EMU values land exactly on target.
$emuEpsilon = 0.00004
$shp.Width  = (5041430 / 12700.0) + $emuEpsilon
$shp.Height = (738623  / 12700.0) + $emuEpsilon

# --- 2. Insert a new first paragraph: "US Core Simple Observation Profile"
$firstPara = $tr.Paragraphs(1, 1)
$firstPara.InsertBefore("US Core Simple Observation Profile`r") | Out-Null

$newPara = $tr.Paragraphs(1, 1)
$newPara.Font.Name = "-apple-system"
$newPara.Font.Bold = $false
$newPara.Font.Italic = $false
$newPara.Font.Color.RGB = 5057303   # 0x4D2B17 little-endian for RGB(23,43,77) == 172B4D

# --- 3. Replace the "Clinical Judgment" paragraph's text with the new
#        multi-run text: "SDC Base Questionnaire/US Core " + "QuestionnaireResponse" + " Profile"
$targetPara = $tr.Paragraphs(3, 1)
$targetPara.Text = "SDC Base Questionnaire/US Core QuestionnaireResponse Profile"

$run1Text = "SDC Base Questionnaire/US Core "
$run2Text = "QuestionnaireResponse"
$run3Text = " Profile"

# Touch each sub-range so the engine emits them as distinct runs.
$r1 = $targetPara.Characters(1, $run1Text.Length)
$r1.Font.Name = "Calibri"

$r2 = $targetPara.Characters($run1Text.Length + 1, $run2Text.Length)
$r2.Font.Name = "Calibri"

$r3 = $targetPara.Characters($run1Text.Length + $run2Text.Length + 1, $run3Text.Length)
$r3.Font.Name = "Calibri"
